# Issue #5: cash & deposit done
#
# The "存款" (deposit) sheet (3rd sheet) gets the same per-row metadata
# treatment that the "土地" (land) and "建物" (building) sheets already
# have: row 1 becomes real column headers (it used to be a stray copy of
# row 2's data) and every data row gains property_category / category /
# date / legislator_name / legislator_id / source_file / index columns
# in G:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Row 1: replace the accidental data values with proper headers and add
# the new header cells G1:M1.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Re-apply the header's bold/bordered style across the whole B1:M1 span
# (values are untouched by a formats-only paste).
$ws.Range("B1").Copy()
$ws.Range("B1:M1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Data rows 2:7 - append the metadata columns. All rows share the same
# property_category/category/date/legislator_name/legislator_id/
# source_file values; "index" mirrors column A.
# ---------------------------------------------------------------------
$indexValues = @{ 2 = 47; 3 = 48; 4 = 49; 5 = 50; 6 = 51; 7 = 52 }

foreach ($r in 2..7) {
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # Force text format before assignment so Excel doesn't reinterpret
    # the "2012-04-30" literal as a date serial number.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2012-04-30"

    $ws.Range("J$r").Value = "張曉風"
    $ws.Range("K$r").Value = 1748
    $ws.Range("L$r").Value = "tmp23631"
    $ws.Range("M$r").Value = $indexValues[$r]

    # Re-apply the data-row style (same style already used by B:F) across
    # G:M now that every value is in place.
    $ws.Range("B$r").Copy()
    $ws.Range("G$r`:M$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
